$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values look numeric (e.g. "313.21", "5.390") but must remain stored
# as literal text, matching the workbook's existing inlineStr/text cells.
# Force text number format before assignment so Excel does not coerce the
# value to a number, then clear the format again so no stray style index is
# left behind (column D cells have no explicit style in the source file).
$dChanges = @{
    'D2' = '27.195.35'
    'D3' = '1.848.15'
    'D5' = '313.21'
    'D7' = '0.4630'
    'D8' = '0.3695'
    'D9' = '0.07268'
    'D10' = '0.8866'
    'D11' = '19.91'
    'D12' = '0.07825'
    'D13' = '1.901.09'
    'D14' = '5.388'
    'D15' = '6.501'
    'D16' = '91.52'
    'D18' = '0.000008853'
    'D20' = '27.221.80'
    'D22' = '5.055'
    'D23' = '2.090.89'
    'D24' = '10.52'
    'D25' = '2.034'
    'D26' = '151.32'
    'D27' = '18.35'
    'D28' = '2.026'
    'D29' = '115.60'
    'D30' = '5.010'
    'D31' = '0.08832'
    'D32' = '3.149'
    'D33' = '0.7825'
    'D34' = '4.509'
    'D36' = '2.707'
    'D37' = '1.101'
    'D38' = '0.01943'
    'D39' = '0.05214'
    'D40' = '2.951'
    'D41' = '7.030'
    'D42' = '0.5039'
    'D43' = '0.1613'
    'D44' = '8.503'
    'D45' = '0.4757'
    'D46' = '10.39'
    'D48' = '102.87'
    'D49' = '1.634'
    'D50' = '0.06196'
    'D51' = '65.50'
}

foreach ($cell in $dChanges.Keys) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $dChanges[$cell]
    $rng.ClearFormats()
}

# Column E values (percentages with surrounding spaces) are never parsed as
# numbers by Excel, so they can be assigned directly as text.
$eChanges = @{
    'E2' = '  +0.06%  '
    'E3' = '  +0.01%  '
    'E4' = '  -0.41%  '
    'E5' = '  -0.16%  '
    'E6' = '  -0.33%  '
    'E7' = '  -0.18%  '
    'E8' = '  -0.20%  '
    'E9' = '  -1.48%  '
    'E10' = '  +0.25%  '
    'E11' = '  -0.14%  '
    'E12' = '  -1.48%  '
    'E13' = '  +3.30%  '
    'E14' = '  +0.17%  '
    'E15' = '  -1.35%  '
    'E16' = '  -0.51%  '
    'E17' = '  -0.39%  '
    'E18' = '  -1.01%  '
    'E19' = '  -0.35%  '
    'E20' = '  +0.07%  '
    'E21' = '  -1.47%  '
    'E22' = '  -1.69%  '
    'E23' = '  -0.57%  '
    'E24' = '  -0.49%  '
    'E25' = '  +9.08%  '
    'E26' = '  -1.11%  '
    'E27' = '  -0.42%  '
    'E28' = '  -2.06%  '
    'E30' = '  -2.41%  '
    'E31' = '  -0.42%  '
    'E32' = '  +5.95%  '
    'E33' = '  +5.40%  '
    'E34' = '  +0.86%  '
    'E35' = '  +0.37%  '
    'E36' = '  +6.33%  '
    'E37' = '  +2.02%  '
    'E38' = '  -0.51%  '
    'E39' = '  -1.19%  '
    'E40' = '  -0.56%  '
    'E41' = '  -1.03%  '
    'E42' = '  -2.66%  '
    'E43' = '  -1.49%  '
    'E44' = '  +2.62%  '
    'E45' = '  -2.19%  '
    'E46' = '  +1.14%  '
    'E47' = '  -0.46%  '
    'E48' = '  -0.01%  '
    'E49' = '  +0.12%  '
    'E50' = '  -0.58%  '
    'E51' = '  -0.05%  '
}

foreach ($cell in $eChanges.Keys) {
    $ws.Range($cell).Value = $eChanges[$cell]
}
